# Apply the "GenerateScript.xlsx" tax_slab_details edit:
#  - Insert a "VALUES" keyword (on its own line) before the opening
#    parenthesis of the generated INSERT statement in column F.
#  - Extend the visible selection to F1:F92 with F92 as the anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tax_slab_details")

# Row 2 is a standalone (non-shared) formula.
$f2 = '="INSERT INTO dbo.tax_slab_details' + "`n" + ' (id, tax_slab_id, from_amount, to_amount, percentage)' + "`n" + 'VALUES ("' + ' & A2 & ", "' + '& B2 & ", " & C2 & ", " & D2 & ", " & E2 & ")"'
$ws.Range("F2").Formula = $f2

# Rows 3:66 share one formula (originally si="0", ref="F3:F66").
# Assigning the formula to the whole range at once keeps it as a single
# shared-formula group with relative references re-based per row.
$f3 = '="INSERT INTO dbo.tax_slab_details' + "`n" + ' (id, tax_slab_id, from_amount, to_amount, percentage)' + "`n" + 'VALUES ("' + ' & A3 & ", "' + '& B3 & ", " & C3 & ", " & D3 & ", " & E3 & ")"'
$ws.Range("F3:F66").Formula = $f3

# Rows 67:91 share another formula (originally si="6", ref="F67:F91").
$f67 = '="INSERT INTO dbo.tax_slab_details' + "`n" + ' (id, tax_slab_id, from_amount, to_amount, percentage)' + "`n" + 'VALUES ("' + ' & A67 & ", "' + '& B67 & ", " & C67 & ", " & D67 & ", " & E67 & ")"'
$ws.Range("F67:F91").Formula = $f67

# F92 is a literal label ("SET IDENTITY_INSERT dbo.tax_slab_details OFF")
# and is intentionally left untouched.

# Recalculate so every cached <v> reflects the new formula text.
$excel.Calculate()

# Match the updated view state: selection spans F1:F92 with F92 active.
$ws.Range("F1:F92").Select()
